$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 402.58334
$ws.Range("I2").Value = 246.83333
$ws.Range("K2").Value = 246.83333
$ws.Range("M2").Value = -133.83333
$ws.Range("H9").Value = 133.55556
$ws.Range("I9").Value = 155.71428
$ws.Range("K9").Value = 155.71428
$ws.Range("M9").Value = 13.28572
$ws.Range("H28").Value = 204.45454
$ws.Range("I28").Value = 121.125
$ws.Range("J28").Value = 426.66666
$ws.Range("K28").Value = 121.125
$ws.Range("L28").Value = 426.66666
$ws.Range("M28").Value = 363.875
$ws.Range("N28").Value = -1396.66666
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").ClearContents()
$ws.Range("H38").Value = 602
$ws.Range("I38").Value = 602
$ws.Range("K38").Value = 1806
$ws.Range("M38").Value = -1434
$ws.Range("H43").Value = 1491
$ws.Range("J43").Value = 1488
$ws.Range("L43").Value = 1488
$ws.Range("N43").Value = -1626
$ws.Range("H53").Value = 1364.125
$ws.Range("I53").Value = 2114.8
$ws.Range("K53").Value = 2114.8
$ws.Range("M53").Value = -1477.8
$ws.Range("H58").Value = 1459
$ws.Range("J58").Value = 3500
$ws.Range("L58").Value = 10500
$ws.Range("N58").Value = -10800
$ws.Range("H116").Value = 8499.177
$ws.Range("I116").Value = 13487.889
$ws.Range("J116").Value = 2886.875
$ws.Range("K116").Value = 13487.889
$ws.Range("L116").Value = 2886.875
$ws.Range("M116").Value = -10045.889
$ws.Range("N116").Value = -9770.875
$ws.Range("H135").Value = 669.1111
$ws.Range("I135").Value = 669.1111
$ws.Range("K135").Value = 6021.9999
$ws.Range("M135").Value = -3486.9999
$ws.Range("H137").Value = 33254.97
$ws.Range("J137").Value = 92517.09
$ws.Range("L137").Value = 277551.27
$ws.Range("N137").Value = -282651.27
$ws.Range("H138").Value = 2422.6516
$ws.Range("J138").Value = 2281.9038
$ws.Range("L138").Value = 6845.7114
$ws.Range("N138").Value = -17125.7114

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1535.0714
$ws.Range("J45").Value = 1721.3334
$ws.Range("L45").Value = 1721.3334
$ws.Range("N45").Value = -2475.3334
$ws.Range("H61").Value = 32708.074
$ws.Range("I61").Value = 41181
$ws.Range("K61").Value = 41181
$ws.Range("M61").Value = -40969
$ws.Range("H122").Value = 1187
$ws.Range("I122").Value = 1201.1428
$ws.Range("K122").Value = 3603.4284
$ws.Range("M122").Value = -1153.4284
$ws.Range("H136").Value = 32708.074
$ws.Range("I136").Value = 41181
$ws.Range("K136").Value = 123543
$ws.Range("M136").Value = -120993

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 612.25
$ws.Range("I107").Value = 530.375
$ws.Range("K107").Value = 530.375
$ws.Range("M107").Value = 1389.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1860.1892
$ws.Range("I31").Value = 1531.72
$ws.Range("K31").Value = 1531.72
$ws.Range("M31").Value = -1236.72
$ws.Range("H34").Value = 1860.1892
$ws.Range("I34").Value = 1531.72
$ws.Range("K34").Value = 1531.72
$ws.Range("M34").Value = -1329.72
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H105").Value = 1076.4166
$ws.Range("I105").Value = 1113.1111
$ws.Range("K105").Value = 1113.1111
$ws.Range("M105").Value = 633.8888999999999
$ws.Range("H107").Value = 520.087
$ws.Range("I107").Value = 520.087
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 520.087
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1399.913
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 3878.75
$ws.Range("I114").Value = 411
$ws.Range("J114").Value = 5034.6665
$ws.Range("K114").Value = 1233
$ws.Range("L114").Value = 15103.9995
$ws.Range("M114").Value = 2021
$ws.Range("N114").Value = -21611.9995
$ws.Range("H116").Value = 62501800
$ws.Range("I116").Value = 1064.25
$ws.Range("K116").Value = 3192.75
$ws.Range("M116").Value = 249.25
$ws.Range("H117").Value = 1450.8
$ws.Range("I117").Value = 555
$ws.Range("J117").Value = 1674.75
$ws.Range("K117").Value = 1665
$ws.Range("L117").Value = 5024.25
$ws.Range("M117").Value = 1777
$ws.Range("N117").Value = -11908.25
$ws.Range("H121").Value = 657.8421
$ws.Range("I121").Value = 490.66666
$ws.Range("K121").Value = 1471.99998
$ws.Range("M121").Value = -161.9999800000001
$ws.Range("H131").Value = 16821.697
$ws.Range("J131").Value = 18440.137
$ws.Range("L131").Value = 55320.41099999999
$ws.Range("N131").Value = -65400.41099999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 3750
$ws.Range("I23").Value = 3750
$ws.Range("K23").Value = 3750
$ws.Range("M23").Value = -3527
$ws.Range("H46").Value = 17421.562
$ws.Range("J46").Value = 17249.666
$ws.Range("L46").Value = 17249.666
$ws.Range("N46").Value = -17561.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2691.5
$ws.Range("I122").Value = 2691.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8074.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5624.5
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 3000
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H107").Value = 907
$ws.Range("I107").Value = 750.8570999999999
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 2252.5713
$ws.Range("L107").Value = 6000
$ws.Range("M107").Value = -332.5712999999996
$ws.Range("N107").Value = -9840
$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -59800
$ws.Range("H126").Value = 1852.1923
$ws.Range("I126").Value = 1613.7894
$ws.Range("J126").Value = 2499.2856
$ws.Range("K126").Value = 4841.3682
$ws.Range("L126").Value = 7497.8568
$ws.Range("M126").Value = -2371.3682
$ws.Range("N126").Value = -12437.8568
$ws.Range("H132").Value = 1220.26
$ws.Range("I132").Value = 1071.1842
$ws.Range("K132").Value = 3213.5526
$ws.Range("M132").Value = -683.5526
